# commit code UI + logic
# - "bill" row repurposed as "order" / "order Detail"
# - new "x" markers added in column F for the first three action rows
# - selection moved to D11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 used to be the "bill" row with an empty second column;
# it becomes the "order" / "order Detail" row.
$ws.Range("C10").Value = "order"
$ws.Range("D10").Value = "order Detail"

# New column F "x" flags for product/category/NSX rows.
$ws.Range("F5").Value = "x"
$ws.Range("F6").Value = "x"
$ws.Range("F7").Value = "x"

# Leave the selection where the author left it.
$ws.Range("D11").Select()
